$d = $word.ActiveDocument
$d.Content.Find.Execute("Michelle Arnetta and Tom Coleman", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Michelle Arnetta and Tom Coleman", 2)
